# Apply scheduled market-data refresh updates to leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 551.5714
$ws.Range("I6").Value = 145.25
$ws.Range("J6").Value = 1093.3334
$ws.Range("K6").Value = 435.75
$ws.Range("L6").Value = 3280.0002
$ws.Range("M6").Value = -323.75
$ws.Range("N6").Value = -3504.0002

$ws.Range("H31").Value = 2333
$ws.Range("I31").Value = 1999.5
$ws.Range("K31").Value = 5998.5
$ws.Range("M31").Value = -5768.5

$ws.Range("H39").Value = 483.73914
$ws.Range("I39").Value = 65.35714
$ws.Range("J39").Value = 1134.5555
$ws.Range("K39").Value = 196.07142
$ws.Range("L39").Value = 3403.6665
$ws.Range("M39").Value = 99.92858000000001
$ws.Range("N39").Value = -3995.6665

$ws.Range("H51").Value = 3464
$ws.Range("I51").Value = 1924.2858
$ws.Range("J51").Value = 4031.2632
$ws.Range("K51").Value = 1924.2858
$ws.Range("L51").Value = 4031.2632
$ws.Range("M51").Value = -1440.2858
$ws.Range("N51").Value = -4999.263199999999

$ws.Range("H138").Value = 2367.344
$ws.Range("I138").Value = 1166.4286
$ws.Range("J138").Value = 3356.3333
$ws.Range("K138").Value = 3499.2858
$ws.Range("L138").Value = 10068.9999
$ws.Range("M138").Value = 1640.7142
$ws.Range("N138").Value = -20348.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1070
$ws.Range("I97").Value = 1070
$ws.Range("K97").Value = 1070
$ws.Range("M97").Value = -574

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2604.2104
$ws.Range("I86").Value = 2698.4614
$ws.Range("J86").Value = 2400
$ws.Range("K86").Value = 2698.4614
$ws.Range("L86").Value = 2400
$ws.Range("M86").Value = -1575.4614
$ws.Range("N86").Value = -4646

$ws.Range("H89").Value = 2604.2104
$ws.Range("I89").Value = 2698.4614
$ws.Range("J89").Value = 2400
$ws.Range("K89").Value = 13492.307
$ws.Range("L89").Value = 12000
$ws.Range("M89").Value = -7876.307000000001
$ws.Range("N89").Value = -23232

$ws.Range("H105").Value = 2411.4285
$ws.Range("I105").Value = 2750
$ws.Range("J105").Value = 2276
$ws.Range("K105").Value = 2750
$ws.Range("L105").Value = 2276
$ws.Range("M105").Value = -1003
$ws.Range("N105").Value = -5770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 2400.3333
$ws.Range("I35").Value = 2446.875
$ws.Range("J35").Value = 2028
$ws.Range("K35").Value = 2446.875
$ws.Range("L35").Value = 2028
$ws.Range("M35").Value = -2152.875
$ws.Range("N35").Value = -2616

$ws.Range("H58").Value = 2318.3513
$ws.Range("I58").Value = 2023.5238
$ws.Range("J58").Value = 2705.3125
$ws.Range("K58").Value = 2023.5238
$ws.Range("L58").Value = 2705.3125
$ws.Range("M58").Value = -1820.5238
$ws.Range("N58").Value = -3111.3125

$ws.Range("H62").Value = 4943.5
$ws.Range("I62").Value = 3601.25
$ws.Range("J62").Value = 7628
$ws.Range("K62").Value = 3601.25
$ws.Range("L62").Value = 7628
$ws.Range("M62").Value = -2977.25
$ws.Range("N62").Value = -8876

$ws.Range("H65").Value = 4943.5
$ws.Range("I65").Value = 3601.25
$ws.Range("J65").Value = 7628
$ws.Range("K65").Value = 18006.25
$ws.Range("L65").Value = 38140
$ws.Range("M65").Value = -14886.25
$ws.Range("N65").Value = -44380

$ws.Range("H132").Value = 2391.9048
$ws.Range("I132").Value = 1884.1034
$ws.Range("K132").Value = 5652.3102
$ws.Range("M132").Value = -3122.3102

$ws.Range("H136").Value = 2318.3513
$ws.Range("I136").Value = 2023.5238
$ws.Range("J136").Value = 2705.3125
$ws.Range("K136").Value = 6070.5714
$ws.Range("L136").Value = 8115.9375
$ws.Range("M136").Value = -3520.5714
$ws.Range("N136").Value = -13215.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1855
$ws.Range("J22").Value = 1855
$ws.Range("L22").Value = 5565
$ws.Range("N22").Value = -5903

$ws.Range("H27").Value = 1855
$ws.Range("J27").Value = 1855
$ws.Range("L27").Value = 5565
$ws.Range("N27").Value = -5769

$ws.Range("H46").Value = 1094
$ws.Range("I46").Value = 162.5
$ws.Range("J46").Value = 2336
$ws.Range("K46").Value = 487.5
$ws.Range("L46").Value = 7008
$ws.Range("M46").Value = -396.5
$ws.Range("N46").Value = -7190

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1999.9166
$ws.Range("I80").Value = 1999.8572
$ws.Range("K80").Value = 1999.8572
$ws.Range("M80").Value = -1001.8572

$ws.Range("H83").Value = 1999.9166
$ws.Range("I83").Value = 1999.8572
$ws.Range("K83").Value = 9999.286
$ws.Range("M83").Value = -5007.286

$ws.Range("H122").Value = 2128.2368
$ws.Range("I122").Value = 1670.1923
$ws.Range("K122").Value = 5010.5769
$ws.Range("M122").Value = -2560.5769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 50302.81
$ws.Range("I7").Value = 79427.16
$ws.Range("J7").Value = 2975.75
$ws.Range("K7").Value = 79427.16
$ws.Range("L7").Value = 2975.75
$ws.Range("M7").Value = -79315.16
$ws.Range("N7").Value = -3199.75

$ws.Range("H22").Value = 880.7646999999999
$ws.Range("I22").Value = 515.1667
$ws.Range("J22").Value = 1080.1818
$ws.Range("K22").Value = 515.1667
$ws.Range("L22").Value = 1080.1818
$ws.Range("M22").Value = -220.1667
$ws.Range("N22").Value = -1670.1818

$ws.Range("H27").Value = 880.7646999999999
$ws.Range("I27").Value = 515.1667
$ws.Range("J27").Value = 1080.1818
$ws.Range("K27").Value = 515.1667
$ws.Range("L27").Value = 1080.1818
$ws.Range("M27").Value = -408.1667
$ws.Range("N27").Value = -1294.1818

$ws.Range("H32").Value = 5106.5
$ws.Range("I32").Value = 213
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 213
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = 104
$ws.Range("N32").Value = -10634

$ws.Range("H46").Value = 57301.777
$ws.Range("I46").Value = 127120.125
$ws.Range("J46").Value = 1447.1
$ws.Range("K46").Value = 127120.125
$ws.Range("L46").Value = 1447.1
$ws.Range("M46").Value = -126932.125
$ws.Range("N46").Value = -1823.1

$ws.Range("H55").Value = 235.3
$ws.Range("I55").Value = 183.46153
$ws.Range("K55").Value = 183.46153
$ws.Range("M55").Value = -10.46153000000001

$ws.Range("H61").Value = 537
$ws.Range("I61").Value = 537
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 537
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -335
$ws.Range("N61").ClearContents()

$ws.Range("H68").Value = 1584.125
$ws.Range("I68").Value = 601.2
$ws.Range("J68").Value = 2030.909
$ws.Range("K68").Value = 601.2
$ws.Range("L68").Value = 2030.909
$ws.Range("M68").Value = 147.8
$ws.Range("N68").Value = -3528.909

$ws.Range("H71").Value = 1584.125
$ws.Range("I71").Value = 601.2
$ws.Range("J71").Value = 2030.909
$ws.Range("K71").Value = 3006
$ws.Range("L71").Value = 10154.545
$ws.Range("M71").Value = 738
$ws.Range("N71").Value = -17642.545

$ws.Range("H93").Value = 15607.429
$ws.Range("I93").Value = 18092
$ws.Range("J93").Value = 700
$ws.Range("K93").Value = 18092
$ws.Range("L93").Value = 700
$ws.Range("M93").Value = -16844
$ws.Range("N93").Value = -3196

$ws.Range("H100").Value = 1950
$ws.Range("I100").Value = 1950
$ws.Range("K100").Value = 1950
$ws.Range("M100").Value = -1409

$ws.Range("H113").Value = 537
$ws.Range("I113").Value = 537
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 537
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1633
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 15878056
$ws.Range("I122").Value = 37038704
$ws.Range("J122").Value = 7570
$ws.Range("K122").Value = 111116112
$ws.Range("L122").Value = 22710
$ws.Range("M122").Value = -111113662
$ws.Range("N122").Value = -27610

$ws.Range("H126").Value = 50302.81
$ws.Range("I126").Value = 79427.16
$ws.Range("J126").Value = 2975.75
$ws.Range("K126").Value = 238281.48
$ws.Range("L126").Value = 8927.25
$ws.Range("M126").Value = -235811.48
$ws.Range("N126").Value = -13867.25

$ws.Range("H132").Value = 11921.482
$ws.Range("J132").Value = 10441.286
$ws.Range("L132").Value = 31323.858
$ws.Range("N132").Value = -36383.858

$ws.Range("H136").Value = 13415818
$ws.Range("I136").Value = 102230.55
$ws.Range("J136").Value = 66670170
$ws.Range("K136").Value = 306691.65
$ws.Range("L136").Value = 200010510
$ws.Range("M136").Value = -304141.65
$ws.Range("N136").Value = -200015610

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 3000
$ws.Range("N62").Value = -4248
$ws.Range("M62").Value = -2376

$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 15000
$ws.Range("N65").Value = -21240
$ws.Range("M65").Value = -11880

$ws.Range("H81").Value = 2241.7144
$ws.Range("I81").Value = 3222.5
$ws.Range("J81").Value = 1849.4
$ws.Range("K81").Value = 6445
$ws.Range("L81").Value = 3698.8
$ws.Range("M81").Value = -5384
$ws.Range("N81").Value = -5820.8

$ws.Range("H84").Value = 2241.7144
$ws.Range("I84").Value = 3222.5
$ws.Range("J84").Value = 1849.4
$ws.Range("K84").Value = 32225
$ws.Range("L84").Value = 18494
$ws.Range("M84").Value = -26921
$ws.Range("N84").Value = -29102

$ws.Range("H96").Value = 1237.1666
$ws.Range("J96").Value = 1496.6666
$ws.Range("L96").Value = 1496.6666
$ws.Range("N96").Value = -4242.6666

$ws.Range("H122").Value = 113866.445
$ws.Range("I122").Value = 252876
$ws.Range("J122").Value = 2658.8
$ws.Range("K122").Value = 758628
$ws.Range("L122").Value = 7976.400000000001
$ws.Range("M122").Value = -756178
$ws.Range("N122").Value = -12876.4

$ws.Range("H126").Value = 94599.60000000001
$ws.Range("I126").Value = 108861.08
$ws.Range("K126").Value = 326583.24
$ws.Range("M126").Value = -324113.24
